$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(327).Insert()

$ws.Cells.Item(327, 1).Value = 3
$ws.Cells.Item(327, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(327, 3).Value = "Coquimbo"
$ws.Cells.Item(327, 4).Value = 44694
$ws.Cells.Item(327, 5).Value = 5
$ws.Cells.Item(327, 6).Value = 100112017
$ws.Cells.Item(327, 7).Value = "Apio"
$ws.Cells.Item(327, 8).Value = "Americana (o)"
$ws.Cells.Item(327, 9).Value = "Primera"
$ws.Cells.Item(327, 10).Value = 230
$ws.Cells.Item(327, 11).Value = 8500
$ws.Cells.Item(327, 12).Value = 9000
$ws.Cells.Item(327, 13).Value = 8739
$ws.Cells.Item(327, 14).Value = "`$/docena de matas"
$ws.Cells.Item(327, 15).Value = "Pan de Az$([char]0xFA)car"
$ws.Cells.Item(327, 16).Value = 1456
$ws.Cells.Item(327, 17).Value = 6
$ws.Cells.Item(327, 18).Value = "Hortaliza"
